$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the header label the same (text is unchanged, only shared-string
# bookkeeping changed upstream, which Excel manages automatically).
$ws.Range("A1").Value = "HK_R_acc_G"

# New values for A2:A50 (results recomputed with new HKlist).
$values = @(
    52.871024734982328,
    52.871024734982328,
    52.915194346289752,
    49.911660777385158,
    50.353356890459366,
    51.32508833922261,
    54.196113074204945,
    54.946996466431095,
    54.284452296819786,
    54.328621908127204,
    59.717314487632514,
    59.540636042402831,
    53.710247349823327,
    54.593639575971729,
    53.975265017667837,
    53.798586572438168,
    55.21201413427562,
    55.079505300353361,
    53.312720848056536,
    53.091872791519435,
    53.445229681978802,
    56.625441696113079,
    61.528268551236756,
    61.130742049469966,
    54.284452296819786,
    54.284452296819786,
    54.593639575971729,
    61.439929328621915,
    60.821554770318023,
    54.372791519434628,
    56.71378091872792,
    52.340989399293292,
    52.561837455830386,
    55.697879858657238,
    56.051236749116605,
    58.61307420494699,
    56.537102473498237,
    56.139575971731446,
    57.243816254416956,
    55.344522968197886,
    55.123674911660778,
    55.521201413427555,
    54.196113074204945,
    54.372791519434628,
    53.533568904593643,
    54.063604240282679,
    52.606007067137803,
    54.107773851590103,
    50.706713780918733
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
